$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New instruments' first-candle rows, appended below the existing table
# (rows 16-27 were blank placeholder rows in the original sheet).
#
# Columns: Name | FirstCandleDate | Connection
# Binance rows use the "dd.mm.yy h:mm:ss" format (numFmtId 166, like the
# existing Binance rows), Bybit rows use "d.m.yyyy h:mm:ss" (numFmtId 165),
# matching the pattern already present in rows 2-15 - except where the
# source data calls for the other format explicitly.

$d16 = Get-Date -Year 2022 -Month 4  -Day 15 -Hour 3  -Minute 30 -Second 0
$d17 = Get-Date -Year 2021 -Month 10 -Day 12 -Hour 3  -Minute 20 -Second 0
$d18 = Get-Date -Year 2020 -Month 1  -Day 17 -Hour 8  -Minute 0  -Second 0
$d19 = Get-Date -Year 2020 -Month 10 -Day 21 -Hour 9  -Minute 29 -Second 59
$d20 = Get-Date -Year 2020 -Month 10 -Day 15 -Hour 8  -Minute 0  -Second 0
$d21 = Get-Date -Year 2021 -Month 10 -Day 11 -Hour 4  -Minute 30 -Second 0
$d22 = Get-Date -Year 2020 -Month 2  -Day 7  -Hour 12 -Minute 10 -Second 0
$d23 = Get-Date -Year 2021 -Month 10 -Day 11 -Hour 3  -Minute 55 -Second 0
$d24 = Get-Date -Year 2020 -Month 9  -Day 18 -Hour 7  -Minute 0  -Second 0
$d25 = Get-Date -Year 2021 -Month 11 -Day 30 -Hour 0  -Minute 0  -Second 0
$d26 = Get-Date -Year 2020 -Month 1  -Day 16 -Hour 8  -Minute 0  -Second 0
$d27 = Get-Date -Year 2021 -Month 6  -Day 29 -Hour 7  -Minute 0  -Second 0

$rows = @(
    @(16, "FTTUSDT",  $d16, "dd.mm.yy h:mm:ss",   "Binance"),
    @(17, "FTTUSDT",  $d17, "d.m.yyyy h:mm:ss",   "Bybit"),
    @(18, "LINKUSDT", $d18, "dd.mm.yy h:mm:ss",   "Binance"),
    @(19, "LINKUSDT", $d19, "d.m.yyyy h:mm:ss",   "Bybit"),
    @(20, "NEARUSDT", $d20, "d.m.yyyy h:mm:ss",   "Binance"),
    @(21, "NEARUSDT", $d21, "d.m.yyyy h:mm:ss",   "Bybit"),
    @(22, "ATOMUSDT", $d22, "dd.mm.yy h:mm:ss",   "Binance"),
    @(23, "ATOMUSDT", $d23, "d.m.yyyy h:mm:ss",   "Bybit"),
    @(24, "UNIUSDT",  $d24, "dd.mm.yy h:mm:ss",   "Binance"),
    @(25, "UNIUSDT",  $d25, "d.m.yyyy h:mm:ss",   "Bybit"),
    @(26, "ETCUSDT",  $d26, "dd.mm.yy h:mm:ss",   "Binance"),
    @(27, "ETCUSDT",  $d27, "dd.mm.yy h:mm:ss",   "Bybit")
)

foreach ($row in $rows) {
    $r = $row[0]

    $nameCell = $ws.Cells.Item($r, 1)
    $nameCell.Value = $row[1]
    $nameCell.Font.ThemeColor = 1

    $dateCell = $ws.Cells.Item($r, 2)
    $dateCell.Value = $row[2]
    $dateCell.Font.ThemeColor = 1
    $dateCell.NumberFormat = $row[3]

    $connCell = $ws.Cells.Item($r, 3)
    $connCell.Value = $row[4]
    $connCell.Font.ThemeColor = 1
}
